# Applies the workbook edit described by the diff:
#  - Shift all "Lookup" strings and "Data" (date/time) values from
#    30.07.2024 to 29.08.2024 (a 30-day shift)
#  - Update the "Prediction" values for rows 29-85 (interval 28-84)
#    to reflect the new forecast numbers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction values for rows 29 through 85 (keyed by worksheet row number)
$predictions = @{
    29 = 0
    30 = 0
    31 = 0.011
    32 = 0.015
    33 = 0.028
    34 = 0.071
    35 = 0.137
    36 = 0.204
    37 = 0.26
    38 = 0.297
    39 = 0.354
    40 = 0.397
    41 = 0.378
    42 = 0.45
    43 = 0.492
    44 = 0.515
    45 = 0.529
    46 = 0.565
    47 = 0.594
    48 = 0.613
    49 = 0.627
    50 = 0.629
    51 = 0.629
    52 = 0.629
    53 = 0.623
    54 = 0.617
    55 = 0.598
    56 = 0.592
    57 = 0.587
    58 = 0.583
    59 = 0.562
    60 = 0.514
    61 = 0.504
    62 = 0.476
    63 = 0.455
    64 = 0.45
    65 = 0.446
    66 = 0.426
    67 = 0.403
    68 = 0.337
    69 = 0.282
    70 = 0.247
    71 = 0.213
    72 = 0.174
    73 = 0.143
    74 = 0.118
    75 = 0.088
    76 = 0.068
    77 = 0.058
    78 = 0.045
    79 = 0.039
    80 = 0.03
    81 = 0.021
    82 = 0.013
    83 = 0
    84 = 0
    85 = 0
}

$lastRow = 96

for ($r = 2; $r -le $lastRow; $r++) {
    # Shift the date/time serial in column A by 30 days (30.07.2024 -> 29.08.2024)
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 30

    # Shift the "Lookup" text in column D from "30.07.2024<n>" to "29.08.2024<n>"
    $cellD = $ws.Cells.Item($r, 4)
    $newLookup = $cellD.Value() -replace '^30\.07\.2024', '29.08.2024'
    $cellD.Value() = $newLookup

    # Update prediction values where changed
    if ($predictions.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value() = $predictions[$r]
    }
}
